# Daily attendance processing - 2026-01-11 10:35:15
# Re-orders the "Recorded By" (column G) list so that any entry referring to
# the automated "System" user (case-insensitive match on "system") is moved
# to the front of the comma-separated list, preserving the relative order of
# the remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

# Column G = "Recorded By" (7th column). Data starts on row 2 (row 1 is the header).
$col = 7

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $current = $cell.Text

    if ([string]::IsNullOrEmpty($current)) {
        continue
    }

    $parts = $current.Split(",")

    $systemItems = @()
    $otherItems = @()

    foreach ($part in $parts) {
        $trimmed = $part.Trim()
        if ($trimmed.Length -eq 0) {
            continue
        }
        if ($trimmed.ToLower().Contains("system")) {
            $systemItems += $trimmed
        } else {
            $otherItems += $trimmed
        }
    }

    $reordered = $systemItems + $otherItems
    $newValue = $reordered -join ", "

    if ($newValue -ne $current) {
        $cell.Value = $newValue
    }
}
